# Updated cryptos list with latest prices / 1h volume change percentages.
# D-column values are written with a leading apostrophe (quote-prefix) so
# Excel stores numeric-looking price strings (e.g. "198.62", "1.00") as text,
# matching the workbook's existing inline-string cell layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''76.322.12'
$ws.Range("E2").Value = '  +0.32%  '
$ws.Range("D3").Value = '''3.039.18'
$ws.Range("E3").Value = '  +3.40%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''198.62'
$ws.Range("E5").Value = '  -2.23%  '
$ws.Range("D6").Value = '''619.72'
$ws.Range("E6").Value = '  +3.35%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -1.46%  '
$ws.Range("E9").Value = '  +4.52%  '
$ws.Range("D10").Value = '''3.037.24'
$ws.Range("E10").Value = '  +3.29%  '
$ws.Range("E11").Value = '  -2.13%  '
$ws.Range("E12").Value = '  -0.79%  '
$ws.Range("D13").Value = '''5.26'
$ws.Range("E13").Value = '  +6.10%  '
$ws.Range("D14").Value = '''3.599.25'
$ws.Range("E14").Value = '  +3.21%  '
$ws.Range("D15").Value = '''28.79'
$ws.Range("E15").Value = '  +1.07%  '
$ws.Range("D16").Value = '''76.327.18'
$ws.Range("E16").Value = '  +0.42%  '
$ws.Range("D17").Value = '''0.0000193'
$ws.Range("E17").Value = '  +0.93%  '
$ws.Range("D18").Value = '''3.039.01'
$ws.Range("E18").Value = '  +3.48%  '
$ws.Range("D19").Value = '''13.48'
$ws.Range("E19").Value = '  +1.07%  '
$ws.Range("D20").Value = '''8.93'
$ws.Range("E20").Value = '  +0.99%  '
$ws.Range("D21").Value = '''378.31'
$ws.Range("E21").Value = '  +0.85%  '
$ws.Range("E22").Value = '  -1.20%  '
$ws.Range("D23").Value = '''4.34'
$ws.Range("E23").Value = '  -0.55%  '
$ws.Range("D24").Value = '''3.198.53'
$ws.Range("E24").Value = '  +3.35%  '
$ws.Range("D25").Value = '''72.89'
$ws.Range("E25").Value = '  +1.39%  '
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("E27").Value = '  -0.91%  '
$ws.Range("D28").Value = '''9.73'
$ws.Range("E28").Value = '  -0.47%  '
$ws.Range("E29").Value = '  -1.12%  '
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("D31").Value = '''8.27'
$ws.Range("E31").Value = '  +3.68%  '
$ws.Range("E32").Value = '  -0.25%  '
$ws.Range("E33").Value = '  +3.56%  '
$ws.Range("D34").Value = '''491.42'
$ws.Range("E34").Value = '  -2.55%  '
$ws.Range("D36").Value = '''20.56'
$ws.Range("E36").Value = '  +0.84%  '
$ws.Range("D37").Value = '''162.27'
$ws.Range("E37").Value = '  -1.75%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = '''0.117'
$ws.Range("E38").Value = '  +3.40%  '
$ws.Range("B39").Value = 'WhiteBITCoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D39").Value = '''20.03'
$ws.Range("E39").Value = '  +1.89%  '
$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").Value = '''0.381'
$ws.Range("E40").Value = '  +1.55%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = '''191.28'
$ws.Range("E41").Value = '  +5.78%  '
$ws.Range("E42").Value = '  -5.04%  '
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").Value = '''0.800'
$ws.Range("E44").Value = '  +20.71%  '
$ws.Range("E45").Value = '  +0.82%  '
$ws.Range("E46").Value = '  +4.21%  '
$ws.Range("D47").Value = '''41.60'
$ws.Range("E47").Value = '  +3.62%  '
$ws.Range("E48").Value = '  -2.22%  '
$ws.Range("D49").Value = '''2.41'
$ws.Range("E49").Value = '  +2.31%  '
$ws.Range("D50").Value = '''0.600'
$ws.Range("E50").Value = '  +2.64%  '
$ws.Range("E51").Value = '  -0.23%  '
